# Add a "Turkey" test-data sheet, cloned from the existing "Spain" sheet
# (same template: headers, merges, styles, number formats), then filled in
# with the Turkey-specific market name and ticket reference.

$wb = $excel.ActiveWorkbook
$spain = $wb.Worksheets.Item("Spain")

# Clone Spain (last sheet) to the end of the workbook - this preserves all
# of the template's styles/merged cells/borders/fills exactly.
$spain.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$turkey = $wb.Worksheets.Item($wb.Worksheets.Count)
$turkey.Name = "Turkey"

# Fill in the Turkey-specific market name and user-story / ticket reference.
$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3291"

# The new sheet doesn't keep Spain's widened column C (back to the sheet's
# standard width), and column D is a touch narrower than Spain's.
$turkey.Columns.Item(3).ColumnWidth = 8.43
$turkey.Columns.Item(4).ColumnWidth = 22.6

# Rows 3 and 5 go back to the sheet's default (auto) height; row 4 (which
# holds the larger-font ticket reference) keeps a slightly taller height.
$turkey.Rows.Item(3).AutoFit()
$turkey.Rows.Item(4).AutoFit()
$turkey.Rows.Item(5).AutoFit()
$turkey.Rows.Item(4).RowHeight = 15.6

# Spain is no longer the active sheet: it keeps its whole table selected,
# while Turkey becomes the active tab with its own selection.
$spain.Activate() | Out-Null
$spain.Range("A1:D10").Select() | Out-Null

$turkey.Activate() | Out-Null
$turkey.Range("G10").Select() | Out-Null
